$wb = $excel.ActiveWorkbook

$handbackFile = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9c2e80cbe1cc21b61c4606e7fe9a78d11d2ed8f/e2e/4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81b84ef7e0e8a28c0cffd1669f671a748f86c30b/e2e/4e26e95f-30d7-4575-b028-0c8920bc4c7f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9c2e80cbe1cc21b61c4606e7fe9a78d11d2ed8f/e2e/4e26e95f-30d7-4575-b028-0c8920bc4c7f.md."

# --- zh-cn sheet: row 5 is the 4e26e95f-30d7-4575-b028-0c8920bc4c7f handback record ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I5").Value = $handbackFile
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $targetUrl, "", "", $handbackFile)

$wsZh.Range("J5").Value = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-10-13 13:25:40"
$wsZh.Range("P5").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet: row 5 is the 4e26e95f-30d7-4575-b028-0c8920bc4c7f handback record ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I5").Value = $handbackFile
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $targetUrl, "", "", $handbackFile)

$wsDe.Range("J5").Value = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.de-de.xlf"
$wsDe.Range("K5").Value = "2016-10-13 13:25:56"
$wsDe.Range("P5").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.14
